# edit.ps1 - applies the "SignUp" sheet addition + scenario updates to the
# login-automation workbook, per commit: "updated the keywordEngine file and
# updated the scenarios and added the firefox driver code"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# ---------------------------------------------------------------------------
# 1) Fix the password value used in the login scenario (row 5, column Value)
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "M@hadevia0"

# ---------------------------------------------------------------------------
# 2) Insert 4 new rows (7-10) before "close browser" row, describing the new
#    post-login verification / account / sign-out steps, and push the
#    existing trailing rows down.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).Resize(4, 1).EntireRow.Insert() | Out-Null

$ws.Range("A7").Value = "verify the home page header "
$ws.Range("B7").Value = "xpath"
$ws.Range("C7").Value = "//i18n-string[text()='User Guide']"
$ws.Range("D7").Value = "isDisplayed"
$ws.Range("E7").Value = "NA"

$ws.Range("A8").Value = "get home page header title"
$ws.Range("B8").Value = "xpath"
$ws.Range("C8").Value = "//title"
$ws.Range("D8").Value = "getText"
$ws.Range("E8").Value = "NA"

$ws.Range("A9").Value = "clickaccount name"
$ws.Range("B9").Value = "className"
$ws.Range("C9").Value = "account-name"
$ws.Range("D9").Value = "click"
$ws.Range("E9").Value = "NA"

$ws.Range("A10").Value = "click signout link"
$ws.Range("B10").Value = "id"
$ws.Range("C10").Value = "signout"
$ws.Range("D10").Value = "click"
$ws.Range("E10").Value = "NA"

$ws.Range("A7:E10").Style = $ws.Range("A6").Style

# ---------------------------------------------------------------------------
# 3) Adjust column widths (A widens to fit new longer text, C widens too)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25
$ws.Columns.Item(3).ColumnWidth = 28

# ---------------------------------------------------------------------------
# 4) Add hyperlink on E5 (mirrors the existing hyperlink style used on E3)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E5"), "https://app.hubspot.com/login") | Out-Null

# ---------------------------------------------------------------------------
# 5) Add the new "SignUp" worksheet after "login"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "SignUp"

# Selection bookkeeping to mirror final state as closely as possible
$ws.Range("K18").Select()
$newSheet.Range("A2").Select()

$wb.Save()
